# The workbook has a "Test" sheet and an "O_TransactionActivity" sheet.
# O_TransactionActivity is the active sheet; update the Source column
# (F3:F24) from the old loan id "LOAN1" to the new instrument id
# "IDHJ-EGNY", and move the active-cell selection to A4.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("O_TransactionActivity")
$ws.Activate()

$ws.Range("F3:F24").Value = "IDHJ-EGNY"

$ws.Range("A4").Select()
